# Apply cell value updates per the source diff.
# Cells whose new value is a plain number-like string (single dot, digits
# only) are force-written as Text so Excel keeps them as strings exactly
# as they were originally (avoids numeric auto-coercion / float drift, and
# the style is reset back to Normal afterwards so no formatting changes
# are introduced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.064.31"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "3.548.79"
$ws.Range("E3").Value = "  +3.78%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.29%  "
$ws.Range("D7").Value = "3.545.26"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("E10").Value = "  +3.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.386"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.51%  "
$ws.Range("D13").Value = "4.155.32"
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000183"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.25%  "
$ws.Range("D16").Value = "3.538.65"
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "64.902.57"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.35%  "
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("E21").Value = "  +6.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.27%  "
$ws.Range("E23").Value = "  +6.70%  "
$ws.Range("D24").Value = "3.690.80"
$ws.Range("E24").Value = "  +3.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000115"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.97%  "
$ws.Range("D32").Value = "3.561.63"
$ws.Range("E33").Value = "  +18.02%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.145"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.36%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.93%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0805"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.12%  "
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +18.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  +10.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.22%  "
$ws.Range("D49").Value = "2.446.39"
$ws.Range("E49").Value = "  +12.94%  "
$ws.Range("E50").Value = "  +7.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "303.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.67%  "
